$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 13 (pushes "Programa resumido:" and
#     everything below it down by one row). Row heights / formatting of the
#     shifted rows travel with them automatically. ---
$ws.Rows("13").Insert()

# Copy column formatting (styles) from an existing B/C pair down into the
# freshly inserted row 13 so the new cells pick up the right style indices
# (col B -> style 2, col C -> style 3) before we stamp in values.
$ws.Range("B10").Copy($ws.Range("B13"))
$ws.Range("C10").Copy($ws.Range("C13"))

# The row insert leaves a phantom styled-but-empty A13 cell (inherited from
# A12's formatting) - row 13 in the target has no A cell at all, so clear it.
$ws.Range("A13").Clear()

# --- Row 10 (Objetivos:): replace the misplaced "Marcos Villela Barcza"
#     text with the actual Portuguese objectives text. ---
$ws.Range("B10").Value = "Conferir aos alunos uma noção prática das aplicações, à escala industrial, de processos e produtos estudados em disciplinas de química orgânica e ainda uma visão global das matérias primas mais importantes na área da indústria química de base orgânica."
$ws.Range("C10").Value = "Conferir aos alunos uma noção prática das aplicações, à escala industrial, de processos e produtos estudados em disciplinas de química orgânica e ainda uma visão global das matérias primas mais importantes na área da indústria química de base orgânica."

# --- New row 13 (Docentes responsáveis: value row): the professor name that
#     used to live under "Objetivos:" now lives here. ---
$ws.Range("B13").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("C13").Value = "1285870 - Marcos Villela Barcza"

# --- Row 14 (Programa resumido:): replace "Semestral" with the actual short
#     syllabus text (Portuguese). ---
$ws.Range("B14").Value = "Petróleo, Gás Natural e Petroquímica; Química Fina; Processos Unitários Orgânicos: Nitração; Esterificação; Alquilação e Acilação; Hidrogenação; Sulfonação/Sulfatação; Oxidação."
$ws.Range("C14").Value = "Petróleo, Gás Natural e Petroquímica; Química Fina; Processos Unitários Orgânicos: Nitração; Esterificação; Alquilação e Acilação; Hidrogenação; Sulfonação/Sulfatação; Oxidação."

# --- Row 16 (Programa:): replace the wrongly duplicated "01/01/2022" value
#     with the real detailed program text (Portuguese). ---
$ws.Range("B16").Value = "Petróleo, Gás Natural e Petroquímica; 2- Química Fina: Características, Química Fina X Química de Base, Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unitários Orgânicos: 3.1- Nitração; 3.2- Esterificação; 3.3- Alquilação e Acilação; 3.4- Hidrogenação; 3.5- Sulfonação/Sulfatação; 3.6- Oxidação."
$ws.Range("C16").Value = "Petróleo, Gás Natural e Petroquímica; 2- Química Fina: Características, Química Fina X Química de Base, Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unitários Orgânicos: 3.1- Nitração; 3.2- Esterificação; 3.3- Alquilação e Acilação; 3.4- Hidrogenação; 3.5- Sulfonação/Sulfatação; 3.6- Oxidação."

# --- Row 19 (Método:): replace the misplaced "Marcos Villela Barcza" text
#     with the teaching-method description. ---
$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."

# --- Row 20 (Critério:): replace the teaching-method text that had landed
#     here with the grading-criteria description. ---
$ws.Range("B20").Value = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."

# --- Row 21 (Norma de recuperação:): replace the grading-criteria text with
#     the recovery-exam rule description. ---
$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."

# --- Row 22 (Bibliografia:): replace the recovery-exam-rule text that had
#     landed here with the actual bibliography text. ---
$ws.Range("B22").Value = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petróleo & Energia, São Paulo, SP: , v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."
$ws.Range("C22").Value = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petróleo & Energia, São Paulo, SP: , v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."

# --- Fix the overlapping column definition: column A's custom width should
#     only apply to column 1, not bleed into column 2 (which has its own
#     width/style definition right after it). ---
$ws.Columns("A").ColumnWidth = 30.7109375
